# CreatePriceList test case added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are entered in alphabetical order so that the shared-string table
# ends up sorted, matching the order rows reference those strings by index.
$ws.Range("A3").Value = "App Settings"
$ws.Range("A1").Value = "Functions.CrmChangeArea.name"
$ws.Range("A5").Value = "Help and Support"
$ws.Range("A2").Value = "Sales"
$ws.Range("A4").Value = "Sales Insights settings"

# Header cell is bold (introduces the second font / cellXf in styles.xml)
$ws.Range("A1").Font.Bold = $true

# Widen column A to fit the header text (target stored width ~29.09 chars;
# the engine quantizes ColumnWidth to 1/6-character steps, so 28.3 is the
# input that lands closest on the nearest achievable stored width)
$ws.Columns.Item(1).ColumnWidth = 28.3

# Leave selection on the cell below the last populated row
$ws.Range("A6").Select() | Out-Null
